$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing column-H number format (style index used by H1:H21)
# down onto the new calibration rows before filling in their values.
$ws.Range("H1").Copy()
$ws.Range("H28:H38").PasteSpecial(-4122)

$values = @(
    0.81766399999999995,
    0.98443800000000004,
    0.81528100000000003,
    0.88751100000000005,
    0.78368800000000005,
    0.76649500000000004,
    0.88601600000000003,
    0.86872000000000005,
    0.86872000000000005,
    1.0224899999999999,
    0.84937700000000005
)

$startRow = 28
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value = $values[$i]
}

# New average offset row
$avgCell = $ws.Cells.Item(44, 8)
$avgCell.Formula = "=AVERAGE(H28:H43)"

# Move the active selection to match the new last-edited cell
$ws.Range("J37").Select()
